$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Apply yellow fill (matches existing style used for A22:A33) to the full new range A22:A44
$ws.Range("A22:A44").Interior.Color = 65535

# Write the new/changed "cung vi" labels in column A first, in the same order they were
# authored, so the workbook's shared-string table is rebuilt in the same sequence.
$ws.Range("A31").Value = "Tử Vi tọa thủ cung Mệnh ở Mùi"
$ws.Range("A32").Value = "Tử Vi tọa thủ cung Mệnh ở Sửu"
$ws.Range("A22").Value = "Tử Vi tọa thủ cung Mệnh ở Tỵ"
$ws.Range("A33").Value = "Tử Vi tọa thủ cung Mệnh ở Ngọ"
$ws.Range("A34").Value = "Tử Vi tọa thủ cung Mệnh ở Dần"
$ws.Range("A35").Value = "Tử Vi tọa thủ cung Mệnh ở Thân"
$ws.Range("A36").Value = "Tử Vi tọa thủ cung Mệnh ở Tuất"
$ws.Range("A23").Value = "Tử Vi tọa thủ cung Mệnh ở Thìn"
$ws.Range("A37").Value = "Tử Vi tọa thủ cung Mệnh ở Tý"
$ws.Range("A38").Value = "Tử Vi tọa thủ cung Mệnh ở Hợi"
$ws.Range("A39").Value = "Tử Vi tọa thủ cung Mệnh ở Mão"
$ws.Range("A40").Value = "Tử Vi tọa thủ cung Mệnh ở Dậu"
$ws.Range("A41").Value = "Tử Vi tọa thủ cung Mệnh gặp Kình Dương"
$ws.Range("A42").Value = "Tử Vi tọa thủ cung Mệnh gặp Đà La"
$ws.Range("A43").Value = "Tử Vi tọa thủ cung Mệnh gặp Địa Không"
$ws.Range("A44").Value = "Tử Vi tọa thủ cung Mệnh gặp Địa Kiếp"

# Now fill in the remaining column A/B/C values for every row in the updated table
# Row 22
$ws.Range("A22").Value = "Tử Vi tọa thủ cung Mệnh ở Tỵ"
$ws.Range("B22").Value = "Bạn là người Thông minh, trung hậu."

# Row 23
$ws.Range("A23").Value = "Tử Vi tọa thủ cung Mệnh ở Thìn"
$ws.Range("B23").Value = "Bạn là người đa mưu, túc trí nhưng vì cái lợi bản thân là phần nhiều."

# Row 24
$ws.Range("A24").Value = "Tử Vi tọa thủ cung Mệnh và hội chiếu các sao Thiên Tướng, Văn Khúc, Văn Xương, Thiên Khôi, Thiên Việt, Tả Phù, Hữu Bật"
$ws.Range("B24").Value = "Bạn là người có uy quyền khiến người khác nể trọng và giúp đỡ. Bản thân ra ngoài gặp nhiều may mắn."

# Row 25
$ws.Range("A25").Value = "Tử Vi tọa thủ cung Mệnh và hội chiếu Thiên Phủ"
$ws.Range("B25").Value = "Bạn có nhiều tiền bạc, của cải."

# Row 26
$ws.Range("A26").Value = "Tử Vi tọa thủ cung Mệnh và gặp Thiên Mã, Lộc Tồn"
$ws.Range("B26").Value = "Độ số quyền lực của bạn được tăng thêm."

# Row 27
$ws.Range("A27").Value = "Tử Vi đồng cung với Thất Sát"
$ws.Range("B27").Value = "Độ số quyền lực của bạn là tuyệt đối."
$ws.Range("C27").Value = "Chế ác được sự tác họa của Hỏa Linh"

# Row 28
$ws.Range("A28").Value = "Tử Vi tọa thủ cung Mệnh và gặp Kình Dương, Đà La"
$ws.Range("B28").Value = "Bạn như vị vua bị vậy hãm."
$ws.Range("C28").Value = "Bị tiểu nhân làm hại."

# Row 29
$ws.Range("A29").Value = "Tử Vi tọa thủ cung Mệnh và gặp Địa Không, Địa Kiếp"
$ws.Range("B29").Value = "Bạn như vị vua bị vây hãm."
$ws.Range("C29").Value = "Bị tiểu nhân làm hại."

# Row 30
$ws.Range("A30").Value = "Tử Vi tọa thủ cung Mệnh và gặp Kình Dương, Đà La, Địa Không, Địa Kiếp"
$ws.Range("B30").Value = "Bạn như vị vua bị vây hãm không lối thoát."
$ws.Range("C30").Value = "Bị tiểu nhân làm hại."

# Row 31
$ws.Range("A31").Value = "Tử Vi tọa thủ cung Mệnh ở Mùi"
$ws.Range("B31").Value = "Bạn là người thông minh, mưu lược, nhưng có phần liều lĩnh."

# Row 32
$ws.Range("A32").Value = "Tử Vi tọa thủ cung Mệnh ở Sửu"
$ws.Range("B32").Value = "Bạn là người thông minh, mưu lược, nhưng có phần liều lĩnh."

# Row 33
$ws.Range("A33").Value = "Tử Vi tọa thủ cung Mệnh ở Ngọ"
$ws.Range("B33").Value = "Bạn là người Thông minh, trung hậu."

# Row 34
$ws.Range("A34").Value = "Tử Vi tọa thủ cung Mệnh ở Dần"
$ws.Range("B34").Value = "Bạn là người Thông minh, trung hậu."

# Row 35
$ws.Range("A35").Value = "Tử Vi tọa thủ cung Mệnh ở Thân"
$ws.Range("B35").Value = "Bạn là người Thông minh, trung hậu."

# Row 36
$ws.Range("A36").Value = "Tử Vi tọa thủ cung Mệnh ở Tuất"
$ws.Range("B36").Value = "Bạn là người đa mưu, túc trí nhưng vì cái lợi bản thân là phần nhiều."

# Row 37
$ws.Range("A37").Value = "Tử Vi tọa thủ cung Mệnh ở Tý"
$ws.Range("B37").Value = "Bạn hơi kém thông minh, nhưng bản tính đôn hậu."
$ws.Range("C37").Value = "Quyền uy kém rực rỡ, khả năng tiêu giảm tai ách bị giảm nhiều."

# Row 38
$ws.Range("A38").Value = "Tử Vi tọa thủ cung Mệnh ở Hợi"
$ws.Range("B38").Value = "Bạn hơi kém thông minh, nhưng bản tính đôn hậu."
$ws.Range("C38").Value = "Quyền uy kém rực rỡ, khả năng tiêu giảm tai ách bị giảm nhiều."

# Row 39
$ws.Range("A39").Value = "Tử Vi tọa thủ cung Mệnh ở Mão"
$ws.Range("B39").Value = "Bạn hơi kém thông minh, nhưng bản tính đôn hậu."
$ws.Range("C39").Value = "Quyền uy kém rực rỡ, khả năng tiêu giảm tai ách bị giảm nhiều."

# Row 40
$ws.Range("A40").Value = "Tử Vi tọa thủ cung Mệnh ở Dậu"
$ws.Range("B40").Value = "Bạn hơi kém thông minh, nhưng bản tính đôn hậu."
$ws.Range("C40").Value = "Quyền uy kém rực rỡ, khả năng tiêu giảm tai ách bị giảm nhiều."

# Row 41
$ws.Range("A41").Value = "Tử Vi tọa thủ cung Mệnh gặp Kình Dương"
$ws.Range("B41").Value = "Bạn như vị vua bị vây hãm."
$ws.Range("C41").Value = "Bị tiểu nhân làm hại."

# Row 42
$ws.Range("A42").Value = "Tử Vi tọa thủ cung Mệnh gặp Đà La"
$ws.Range("B42").Value = "Bạn như vị vua bị vây hãm."
$ws.Range("C42").Value = "Bị tiểu nhân làm hại."

# Row 43
$ws.Range("A43").Value = "Tử Vi tọa thủ cung Mệnh gặp Địa Không"
$ws.Range("B43").Value = "Bạn như vị vua bị vây hãm."
$ws.Range("C43").Value = "Bị tiểu nhân làm hại."

# Row 44
$ws.Range("A44").Value = "Tử Vi tọa thủ cung Mệnh gặp Địa Kiếp"
$ws.Range("B44").Value = "Bạn như vị vua bị vây hãm."
$ws.Range("C44").Value = "Bị tiểu nhân làm hại."

# Clear cells that held content in the old layout but are blank in the updated one
$ws.Range("C25").ClearContents()
$ws.Range("C31").ClearContents()
$ws.Range("C32").ClearContents()
$ws.Range("C33").ClearContents()

# Update sheet view: scroll so row 10 is at top, and set the active selection to J29
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J29").Select()
